# Update the cached date shown in the small "Rectangle 6" date-field shape
# on the slide master (e.g. from "19.06.2024" to "06.07.2024").
$p = $ppt.ActivePresentation

$master = $p.SlideMaster
$dateShape = $master.Shapes.Item("Rectangle 6")
$dateShape.TextFrame.TextRange.Text = "06.07.2024"

# Update the presenter/date line on slide 1 ("19.06.2024, Daniel Krämer"
# -> "10.07.2024, Daniel Krämer"), splitting the old date prefix into its
# own run so that only the date portion is freshly (re)typed while the
# ", Daniel Krämer" suffix remains the original run.
$slide1 = $p.Slides.Item(1)
$infoShape = $slide1.Shapes.Item("Rectangle 3")
$tr = $infoShape.TextFrame.TextRange
$oldDateLen = ("19.06.2024").Length
$dateRange = $tr.Characters(1, $oldDateLen)
$dateRange.Text = "10.07.2024"
